$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(112, 8).Value = 875071.9399999999  # H112: 989084.2 -> 875071.9399999999
$ws.Cells.Item(112, 10).Value = 989164  # J112: 1137391.9 -> 989164
$ws.Cells.Item(112, 12).Value = 2967492  # L112: 3412175.7 -> 2967492
$ws.Cells.Item(112, 14).Value = -2969708  # N112: -3414391.7 -> -2969708

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(10, 8).Value = 8751.25  # H10: 4500 -> 8751.25
$ws.Cells.Item(10, 9).Value = 5000  # I10: 3500 -> 5000
$ws.Cells.Item(10, 10).Value = 10001.667  # J10: 6500 -> 10001.667
$ws.Cells.Item(10, 11).Value = 5000  # K10: 3500 -> 5000
$ws.Cells.Item(10, 12).Value = 10001.667  # L10: 6500 -> 10001.667
$ws.Cells.Item(10, 13).Value = -4830  # M10: -3330 -> -4830
$ws.Cells.Item(10, 14).Value = -10341.667  # N10: -6840 -> -10341.667

$ws.Cells.Item(11, 8).Value = 1700  # H11: 0 -> 1700
$ws.Cells.Item(11, 9).Value = 50  # I11: 0 -> 50
$ws.Cells.Item(11, 10).Value = 5000  # J11: 0 -> 5000
$ws.Cells.Item(11, 11).Value = 50  # K11: 0 -> 50
$ws.Cells.Item(11, 12).Value = 5000  # L11: 0 -> 5000
$ws.Cells.Item(11, 13).Value = 94  # M11: None -> 94
$ws.Cells.Item(11, 14).Value = -5288  # N11: None -> -5288

$ws.Cells.Item(13, 8).Value = 10000  # H13: 0 -> 10000
$ws.Cells.Item(13, 10).Value = 10000  # J13: 0 -> 10000
$ws.Cells.Item(13, 12).Value = 10000  # L13: 0 -> 10000
$ws.Cells.Item(13, 14).Value = -10288  # N13: None -> -10288

$ws.Cells.Item(37, 8).Value = 10190.833  # H37: 7684.3335 -> 10190.833
$ws.Cells.Item(37, 9).Value = 2722.6667  # I37: 2542 -> 2722.6667
$ws.Cells.Item(37, 10).Value = 17659  # J37: 17969 -> 17659
$ws.Cells.Item(37, 11).Value = 2722.6667  # K37: 2542 -> 2722.6667
$ws.Cells.Item(37, 12).Value = 17659  # L37: 17969 -> 17659
$ws.Cells.Item(37, 13).Value = -2449.6667  # M37: -2269 -> -2449.6667
$ws.Cells.Item(37, 14).Value = -18205  # N37: -18515 -> -18205

$ws.Cells.Item(44, 8).Value = 19898  # H44: 19898.111 -> 19898
$ws.Cells.Item(44, 10).Value = 19898  # J44: 19898.111 -> 19898
$ws.Cells.Item(44, 12).Value = 19898  # L44: 19898.111 -> 19898
$ws.Cells.Item(44, 14).Value = -20874  # N44: -20874.111 -> -20874

$ws.Cells.Item(55, 8).Value = 0  # H55: 24699.5 -> 0
$ws.Cells.Item(55, 10).Value = 0  # J55: 24699.5 -> 0
$ws.Cells.Item(55, 12).Value = 0  # L55: 24699.5 -> 0
$ws.Cells.Item(55, 14).ClearContents()  # N55: -25329.5 -> (removed)

$ws.Cells.Item(63, 8).Value = 2243.5715  # H63: 2987.2666 -> 2243.5715
$ws.Cells.Item(63, 9).Value = 2243.5715  # I63: 2255.3635 -> 2243.5715
$ws.Cells.Item(63, 10).Value = 0  # J63: 5000 -> 0
$ws.Cells.Item(63, 11).Value = 2243.5715  # K63: 2255.3635 -> 2243.5715
$ws.Cells.Item(63, 12).Value = 0  # L63: 5000 -> 0
$ws.Cells.Item(63, 13).Value = -1557.5715  # M63: -1569.3635 -> -1557.5715

$ws.Cells.Item(66, 8).Value = 2243.5715  # H66: 2987.2666 -> 2243.5715
$ws.Cells.Item(66, 9).Value = 2243.5715  # I66: 2255.3635 -> 2243.5715
$ws.Cells.Item(66, 10).Value = 0  # J66: 5000 -> 0
$ws.Cells.Item(66, 11).Value = 11217.8575  # K66: 11276.8175 -> 11217.8575
$ws.Cells.Item(66, 12).Value = 0  # L66: 25000 -> 0
$ws.Cells.Item(66, 13).Value = -7785.8575  # M66: -7844.817499999999 -> -7785.8575

$ws.Cells.Item(74, 8).Value = 438306.6  # H74: 315214.34 -> 438306.6
$ws.Cells.Item(74, 9).Value = 4108.625  # I74: 2943.5417 -> 4108.625
$ws.Cells.Item(74, 10).Value = 1430759.1  # J74: 1252026.8 -> 1430759.1
$ws.Cells.Item(74, 11).Value = 4108.625  # K74: 2943.5417 -> 4108.625
$ws.Cells.Item(74, 12).Value = 1430759.1  # L74: 1252026.8 -> 1430759.1
$ws.Cells.Item(74, 13).Value = -3234.625  # M74: -2069.5417 -> -3234.625
$ws.Cells.Item(74, 14).Value = -1432507.1  # N74: -1253774.8 -> -1432507.1

$ws.Cells.Item(77, 8).Value = 438306.6  # H77: 315214.34 -> 438306.6
$ws.Cells.Item(77, 9).Value = 4108.625  # I77: 2943.5417 -> 4108.625
$ws.Cells.Item(77, 10).Value = 1430759.1  # J77: 1252026.8 -> 1430759.1
$ws.Cells.Item(77, 11).Value = 20543.125  # K77: 14717.7085 -> 20543.125
$ws.Cells.Item(77, 12).Value = 7153795.5  # L77: 6260134 -> 7153795.5
$ws.Cells.Item(77, 13).Value = -16175.125  # M77: -10349.7085 -> -16175.125
$ws.Cells.Item(77, 14).Value = -7162531.5  # N77: -6268870 -> -7162531.5

$ws.Cells.Item(80, 8).Value = 25577.777  # H80: 28000 -> 25577.777
$ws.Cells.Item(80, 10).Value = 25577.777  # J80: 28000 -> 25577.777
$ws.Cells.Item(80, 12).Value = 25577.777  # L80: 28000 -> 25577.777
$ws.Cells.Item(80, 14).Value = -27573.777  # N80: -29996 -> -27573.777

$ws.Cells.Item(83, 8).Value = 25577.777  # H83: 28000 -> 25577.777
$ws.Cells.Item(83, 10).Value = 25577.777  # J83: 28000 -> 25577.777
$ws.Cells.Item(83, 12).Value = 76733.33099999999  # L83: 84000 -> 76733.33099999999
$ws.Cells.Item(83, 14).Value = -86717.33099999999  # N83: -93984 -> -86717.33099999999

$ws.Cells.Item(127, 8).Value = 30000  # H127: 29999.273 -> 30000
$ws.Cells.Item(127, 10).Value = 30000  # J127: 29999.273 -> 30000
$ws.Cells.Item(127, 12).Value = 30000  # L127: 29999.273 -> 30000
$ws.Cells.Item(127, 14).Value = -39920  # N127: -39919.273 -> -39920

$ws.Cells.Item(132, 8).Value = 34798.312  # H132: 33750.242 -> 34798.312
$ws.Cells.Item(132, 9).Value = 45360.332  # I132: 41910 -> 45360.332
$ws.Cells.Item(132, 10).Value = 3112.25  # J132: 3442.5715 -> 3112.25
$ws.Cells.Item(132, 11).Value = 136080.996  # K132: 125730 -> 136080.996
$ws.Cells.Item(132, 12).Value = 9336.75  # L132: 10327.7145 -> 9336.75
$ws.Cells.Item(132, 13).Value = -133550.996  # M132: -123200 -> -133550.996
$ws.Cells.Item(132, 14).Value = -14396.75  # N132: -15387.7145 -> -14396.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(26, 8).Value = 16000  # H26: 8760 -> 16000
$ws.Cells.Item(26, 9).Value = 9000  # I26: 8760 -> 9000
$ws.Cells.Item(26, 10).Value = 30000  # J26: 0 -> 30000
$ws.Cells.Item(26, 11).Value = 9000  # K26: 8760 -> 9000
$ws.Cells.Item(26, 12).Value = 30000  # L26: 0 -> 30000
$ws.Cells.Item(26, 13).Value = -8708  # M26: -8468 -> -8708
$ws.Cells.Item(26, 14).Value = -30584  # N26: None -> -30584

$ws.Cells.Item(96, 8).Value = 11203  # H96: 13225.417 -> 11203
$ws.Cells.Item(96, 9).Value = 7403.857  # I96: 9522.777 -> 7403.857
$ws.Cells.Item(96, 10).Value = 24500  # J96: 24333.334 -> 24500
$ws.Cells.Item(96, 11).Value = 7403.857  # K96: 9522.777 -> 7403.857
$ws.Cells.Item(96, 12).Value = 24500  # L96: 24333.334 -> 24500
$ws.Cells.Item(96, 13).Value = -4657.857  # M96: -6776.777 -> -4657.857
$ws.Cells.Item(96, 14).Value = -29992  # N96: -29825.334 -> -29992

$ws.Cells.Item(105, 8).Value = 1871.2609  # H105: 2514.5454 -> 1871.2609
$ws.Cells.Item(105, 9).Value = 1511.2667  # I105: 2010 -> 1511.2667
$ws.Cells.Item(105, 10).Value = 2546.25  # J105: 3120 -> 2546.25
$ws.Cells.Item(105, 11).Value = 1511.2667  # K105: 2010 -> 1511.2667
$ws.Cells.Item(105, 12).Value = 2546.25  # L105: 3120 -> 2546.25
$ws.Cells.Item(105, 13).Value = 235.7333000000001  # M105: -263 -> 235.7333000000001
$ws.Cells.Item(105, 14).Value = -6040.25  # N105: -6614 -> -6040.25

$ws.Cells.Item(130, 8).Value = 29779.666  # H130: 27733.334 -> 29779.666
$ws.Cells.Item(130, 10).Value = 29779.666  # J130: 27733.334 -> 29779.666
$ws.Cells.Item(130, 12).Value = 29779.666  # L130: 27733.334 -> 29779.666
$ws.Cells.Item(130, 14).Value = -39819.666  # N130: -37773.334 -> -39819.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(20, 8).Value = 50390  # H20: 0 -> 50390
$ws.Cells.Item(20, 10).Value = 50390  # J20: 0 -> 50390
$ws.Cells.Item(20, 12).Value = 50390  # L20: 0 -> 50390
$ws.Cells.Item(20, 14).Value = -50862  # N20: None -> -50862

$ws.Cells.Item(30, 8).Value = 50390  # H30: 0 -> 50390
$ws.Cells.Item(30, 10).Value = 50390  # J30: 0 -> 50390
$ws.Cells.Item(30, 12).Value = 50390  # L30: 0 -> 50390
$ws.Cells.Item(30, 14).Value = -50572  # N30: None -> -50572

$ws.Cells.Item(105, 8).Value = 949.75  # H105: 797.75 -> 949.75
$ws.Cells.Item(105, 9).Value = 899  # I105: 380 -> 899
$ws.Cells.Item(105, 10).Value = 966.6667  # J105: 937 -> 966.6667
$ws.Cells.Item(105, 11).Value = 899  # K105: 380 -> 899
$ws.Cells.Item(105, 12).Value = 966.6667  # L105: 937 -> 966.6667
$ws.Cells.Item(105, 13).Value = 848  # M105: 1367 -> 848
$ws.Cells.Item(105, 14).Value = -4460.6667  # N105: -4431 -> -4460.6667

$ws.Cells.Item(122, 8).Value = 1167.6296  # H122: 1283.1904 -> 1167.6296
$ws.Cells.Item(122, 9).Value = 837.5  # I122: 901.2222 -> 837.5
$ws.Cells.Item(122, 10).Value = 1431.7333  # J122: 1569.6666 -> 1431.7333
$ws.Cells.Item(122, 11).Value = 2512.5  # K122: 2703.6666 -> 2512.5
$ws.Cells.Item(122, 12).Value = 4295.199900000001  # L122: 4708.9998 -> 4295.199900000001
$ws.Cells.Item(122, 13).Value = -62.5  # M122: -253.6666 -> -62.5
$ws.Cells.Item(122, 14).Value = -9195.1999  # N122: -9608.9998 -> -9195.1999

$ws.Cells.Item(128, 8).Value = 50390  # H128: 0 -> 50390
$ws.Cells.Item(128, 10).Value = 50390  # J128: 0 -> 50390
$ws.Cells.Item(128, 12).Value = 50390  # L128: 0 -> 50390
$ws.Cells.Item(128, 14).Value = -60350  # N128: None -> -60350

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 3559.2727  # H113: 4935.7144 -> 3559.2727
$ws.Cells.Item(113, 9).Value = 599.3333  # I113: 380 -> 599.3333
$ws.Cells.Item(113, 10).Value = 4669.25  # J113: 5695 -> 4669.25
$ws.Cells.Item(113, 11).Value = 1797.9999  # K113: 1140 -> 1797.9999
$ws.Cells.Item(113, 12).Value = 14007.75  # L113: 17085 -> 14007.75
$ws.Cells.Item(113, 13).Value = 372.0001  # M113: 1030 -> 372.0001
$ws.Cells.Item(113, 14).Value = -18347.75  # N113: -21425 -> -18347.75

$ws.Cells.Item(132, 8).Value = 1020.5769  # H132: 1099 -> 1020.5769
$ws.Cells.Item(132, 9).Value = 510.9375  # I132: 527.6667 -> 510.9375
$ws.Cells.Item(132, 10).Value = 1836  # J132: 1813.1666 -> 1836
$ws.Cells.Item(132, 11).Value = 4598.4375  # K132: 4749.0003 -> 4598.4375
$ws.Cells.Item(132, 12).Value = 16524  # L132: 16318.4994 -> 16524
$ws.Cells.Item(132, 13).Value = -2068.4375  # M132: -2219.0003 -> -2068.4375
$ws.Cells.Item(132, 14).Value = -21584  # N132: -21378.4994 -> -21584

$ws.Cells.Item(138, 8).Value = 1413.4706  # H138: 1908 -> 1413.4706
$ws.Cells.Item(138, 9).Value = 1621.5  # I138: 1102.6364 -> 1621.5
$ws.Cells.Item(138, 10).Value = 1300  # J138: 2892.3333 -> 1300
$ws.Cells.Item(138, 11).Value = 4864.5  # K138: 3307.9092 -> 4864.5
$ws.Cells.Item(138, 12).Value = 3900  # L138: 8676.999899999999 -> 3900
$ws.Cells.Item(138, 13).Value = 275.5  # M138: 1832.0908 -> 275.5
$ws.Cells.Item(138, 14).Value = -14180  # N138: -18956.9999 -> -14180

$ws.Cells.Item(140, 8).Value = 1954.762  # H140: 1913.1818 -> 1954.762
$ws.Cells.Item(140, 9).Value = 1208.8235  # I140: 1199.4445 -> 1208.8235
$ws.Cells.Item(140, 11).Value = 3626.4705  # K140: 3598.3335 -> 3626.4705
$ws.Cells.Item(140, 13).Value = 1553.5295  # M140: 1581.6665 -> 1553.5295

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1967.7241  # H7: 2279.0908 -> 1967.7241
$ws.Cells.Item(7, 9).Value = 1307.4546  # I7: 1685.5555 -> 1307.4546
$ws.Cells.Item(7, 10).Value = 4042.8572  # J7: 4950 -> 4042.8572
$ws.Cells.Item(7, 11).Value = 1307.4546  # K7: 1685.5555 -> 1307.4546
$ws.Cells.Item(7, 12).Value = 4042.8572  # L7: 4950 -> 4042.8572
$ws.Cells.Item(7, 13).Value = -1195.4546  # M7: -1573.5555 -> -1195.4546
$ws.Cells.Item(7, 14).Value = -4266.8572  # N7: -5174 -> -4266.8572

$ws.Cells.Item(16, 8).Value = 0  # H16: 1123.75 -> 0
$ws.Cells.Item(16, 9).Value = 0  # I16: 1123.75 -> 0
$ws.Cells.Item(16, 11).Value = 0  # K16: 1123.75 -> 0
$ws.Cells.Item(16, 13).ClearContents()  # M16: -953.75 -> (removed)

$ws.Cells.Item(46, 8).Value = 1530.3636  # H46: 1473.1538 -> 1530.3636
$ws.Cells.Item(46, 9).Value = 2072  # I46: 2950.25 -> 2072
$ws.Cells.Item(46, 10).Value = 880.4  # J46: 816.6667 -> 880.4
$ws.Cells.Item(46, 11).Value = 2072  # K46: 2950.25 -> 2072
$ws.Cells.Item(46, 12).Value = 880.4  # L46: 816.6667 -> 880.4
$ws.Cells.Item(46, 13).Value = -1884  # M46: -2762.25 -> -1884
$ws.Cells.Item(46, 14).Value = -1256.4  # N46: -1192.6667 -> -1256.4

$ws.Cells.Item(125, 8).Value = 19780.875  # H125: 17828.572 -> 19780.875
$ws.Cells.Item(125, 10).Value = 34999  # J125: 60900 -> 34999
$ws.Cells.Item(125, 12).Value = 34999  # L125: 60900 -> 34999
$ws.Cells.Item(125, 14).Value = -44839  # N125: -70740 -> -44839

$ws.Cells.Item(126, 8).Value = 1967.7241  # H126: 2279.0908 -> 1967.7241
$ws.Cells.Item(126, 9).Value = 1307.4546  # I126: 1685.5555 -> 1307.4546
$ws.Cells.Item(126, 10).Value = 4042.8572  # J126: 4950 -> 4042.8572
$ws.Cells.Item(126, 11).Value = 3922.3638  # K126: 5056.666499999999 -> 3922.3638
$ws.Cells.Item(126, 12).Value = 12128.5716  # L126: 14850 -> 12128.5716
$ws.Cells.Item(126, 13).Value = -1452.3638  # M126: -2586.666499999999 -> -1452.3638
$ws.Cells.Item(126, 14).Value = -17068.5716  # N126: -19790 -> -17068.5716

$ws.Cells.Item(127, 8).Value = 26072  # H127: 21849.285 -> 26072
$ws.Cells.Item(127, 10).Value = 26072  # J127: 21849.285 -> 26072
$ws.Cells.Item(127, 12).Value = 26072  # L127: 21849.285 -> 26072
$ws.Cells.Item(127, 14).Value = -35992  # N127: -31769.285 -> -35992

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 9167.916999999999  # H2: 9201.6 -> 9167.916999999999
$ws.Cells.Item(2, 9).Value = 0  # I2: 6000 -> 0
$ws.Cells.Item(2, 10).Value = 9167.916999999999  # J2: 9430.286 -> 9167.916999999999
$ws.Cells.Item(2, 11).Value = 0  # K2: 6000 -> 0
$ws.Cells.Item(2, 12).Value = 9167.916999999999  # L2: 9430.286 -> 9167.916999999999
$ws.Cells.Item(2, 13).ClearContents()  # M2: -5888 -> (removed)
$ws.Cells.Item(2, 14).Value = -9391.916999999999  # N2: -9654.286 -> -9391.916999999999

$ws.Cells.Item(128, 8).Value = 30000  # H128: 31800 -> 30000
$ws.Cells.Item(128, 10).Value = 30000  # J128: 31800 -> 30000
$ws.Cells.Item(128, 12).Value = 30000  # L128: 31800 -> 30000
$ws.Cells.Item(128, 14).Value = -39960  # N128: -41760 -> -39960
